# Rename the inline pictures in the headers/footers.
#
#   header1.xml (first-page header) : BTec_Logo-Orange   image1.jpg -> image2.jpg
#   footer1.xml (first-page footer) : PearsonLogo         image2.png -> image1.png
#   footer2.xml (default footer)    : PearsonLogo         image2.png -> image1.png
#
# InlineShape has no settable .Name in the Word object model, so the
# well-known trick is used: convert the inline picture to a floating
# shape (which exposes .Name), rename it, then convert it back to an
# inline picture so the layout/anchor is unchanged.
$d = $word.ActiveDocument

function Rename-InlinePicture($range, $index, $newName) {
    $inlineShape = $range.InlineShapes.Item($index)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

$section = $d.Sections.Item(1)

# First-page header (header1.xml) holds the BTEC logo.
$firstHeader = $section.Headers.Item(2)
if ($firstHeader.Exists -and $firstHeader.Range.InlineShapes.Count -ge 1) {
    Rename-InlinePicture $firstHeader.Range 1 "image2.jpg"
}

# First-page footer (footer1.xml, docPr id="3") holds the Pearson logo.
$firstFooter = $section.Footers.Item(2)
if ($firstFooter.Exists -and $firstFooter.Range.InlineShapes.Count -ge 1) {
    Rename-InlinePicture $firstFooter.Range 1 "image1.png"
}

# Default footer (footer2.xml, docPr id="2") holds the Pearson logo.
$defaultFooter = $section.Footers.Item(1)
if ($defaultFooter.Exists -and $defaultFooter.Range.InlineShapes.Count -ge 1) {
    Rename-InlinePicture $defaultFooter.Range 1 "image1.png"
}
